# Add season-record columns (Wins, Losses, Ties) to the DET_2018 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new headers with same style as existing headers ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting (bold, centered, bordered) from an existing header cell
# onto the new header cells so they match the rest of row 1 (style index 1).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-51): every team/player row gets the same season record ---
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 64
    $ws.Cells.Item($row, 31).Value = 98
    $ws.Cells.Item($row, 32).Value = 0
}
